$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder two pairs of countries in the list (labels swap, data stays keyed to row) ---
# Chipre now appears before Letonia
$ws.Range("A92").Value = "Republica de Chipre"
$ws.Range("A93").Value = "Letonia"

# Sri Lanka now appears before Georgia
$ws.Range("A106").Value = "Sri Lanka"
$ws.Range("A107").Value = "Georgia"

function Set-RowValues {
    param($sheet, $row, $values)
    $col = 2
    foreach ($v in $values) {
        $sheet.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

# --- Updated case numbers (new day's figures) ---

# Estados Unidos (row 4)
Set-RowValues $ws 4 @(963747, 3096, 118336, 791047, 15110, 108, 54364)

# Alemania (row 8)
Set-RowValues $ws 8 @(157026, 513, 109800, 41346, 2570, 3, 5880)

# Brasil (row 14)
Set-RowValues $ws 14 @(59479, 283, 29160, 26257, 8318, 17, 4062)

# Republica Dominicana (row 47)
Set-RowValues $ws 47 @(6135, 209, 910, 4947, 144, 5, 278)

# Moldavia (row 58)
Set-RowValues $ws 58 @(3408, 104, 895, 2417, 212, 2, 96)

# Grecia (row 64)
Set-RowValues $ws 64 @(2517, 11, 577, 1806, 46, 4, 134)

# Row 92 (now Republica de Chipre)
Set-RowValues $ws 92 @(817, 7, 148, 655, 15, 0, 14)

# Row 93 (now Letonia)
Set-RowValues $ws 93 @(812, 8, 267, 533, 6, 0, 12)

# Row 106 (now Sri Lanka)
Set-RowValues $ws 106 @(485, 33, 120, 358, 2, 0, 7)

# Row 107 (now Georgia)
Set-RowValues $ws 107 @(485, 29, 139, 340, 6, 1, 6)
